$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.622.50'
$ws.Range("E2").Value = '  +1.09%  '
$ws.Range("D3").Value = '3.393.01'
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '577.39'
$ws.Range("E5").Value = '  +1.24%  '
$ws.Range("D6").Value = '140.34'
$ws.Range("E6").Value = '  -0.94%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '0.474'
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("D9").Value = '7.67'
$ws.Range("E9").Value = '  +2.08%  '
$ws.Range("E10").Value = '  -0.97%  '
$ws.Range("E11").Value = '  -1.83%  '
$ws.Range("D12").Value = '3.971.99'
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("D13").Value = '28.67'
$ws.Range("E13").Value = '  +0.64%  '
$ws.Range("E14").Value = '  +0.74%  '
$ws.Range("D15").Value = '3.388.27'
$ws.Range("E15").Value = '  -0.07%  '
$ws.Range("E16").Value = '  -0.62%  '
$ws.Range("D17").Value = '61.649.44'
$ws.Range("E17").Value = '  +1.09%  '
$ws.Range("D18").Value = '6.15'
$ws.Range("E18").Value = '  -1.18%  '
$ws.Range("D19").Value = '13.65'
$ws.Range("E19").Value = '  -2.38%  '
$ws.Range("D20").Value = '8.96'
$ws.Range("E20").Value = '  -0.50%  '
$ws.Range("D21").Value = '389.83'
$ws.Range("E21").Value = '  +1.31%  '
$ws.Range("D22").Value = '75.38'
$ws.Range("E22").Value = '  +1.79%  '
$ws.Range("E23").Value = '  -0.51%  '
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("D25").Value = '0.0000112'
$ws.Range("E25").Value = '  -3.92%  '
$ws.Range("E26").Value = '  +6.36%  '
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("D28").Value = '7.24'
$ws.Range("E28").Value = '  -2.34%  '
$ws.Range("D29").Value = '8.05'
$ws.Range("E29").Value = '  +0.58%  '
$ws.Range("E30").Value = '  -0.20%  '
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("E32").Value = '  -4.05%  '
$ws.Range("D33").Value = '23.41'
$ws.Range("E33").Value = '  -1.11%  '
$ws.Range("D34").Value = '6.94'
$ws.Range("E34").Value = '  -1.12%  '
$ws.Range("D35").Value = '167.67'
$ws.Range("E35").Value = '  +1.17%  '
$ws.Range("D36").Value = '5.02'
$ws.Range("E36").Value = '  +0.45%  '
$ws.Range("D37").Value = '3.429.16'
$ws.Range("E37").Value = '  +0.20%  '
$ws.Range("E38").Value = '  -1.10%  '
$ws.Range("E39").Value = '  -1.15%  '
$ws.Range("D40").Value = '25.71'
$ws.Range("E40").Value = '  -9.01%  '
$ws.Range("E41").Value = '  -0.18%  '
$ws.Range("D42").Value = '4.43'
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("E43").Value = '  -0.26%  '
$ws.Range("E44").Value = '  -0.55%  '
$ws.Range("D45").Value = '2.460.65'
$ws.Range("E45").Value = '  -1.13%  '
$ws.Range("E46").Value = '  -1.86%  '
$ws.Range("D47").Value = '22.55'
$ws.Range("E47").Value = '  -3.78%  '
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("D49").Value = '0.0262'
$ws.Range("E49").Value = '  -3.63%  '
$ws.Range("D50").Value = '2.02'
$ws.Range("E50").Value = '  -2.62%  '
$ws.Range("E51").Value = '  -1.67%  '
